$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("I5").Value = 15
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 2
$ws.Range("M5").Value = 2020
